$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "regular_season"

# Update the header cells: "Division" -> "Away Division" (E1) / "Home Division" (G1)
$ws.Range("E1").Value = "Away Division"
$ws.Range("G1").Value = "Home Division"

# Move the selection / scrolled view to H7 (also clears any stale scroll position)
$ws.Range("H7").Select()

# The longer header text widens columns E and G (Excel auto-fit behavior)
$ws.Columns.Item(5).ColumnWidth = 11.0
$ws.Columns.Item(7).ColumnWidth = 11.666666666666666

Write-Output "done"
